# Update the "dSF" column (F) values to reflect the repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = 5
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 6
$ws.Range("F11").Value = -7
$ws.Range("F12").Value = -7
$ws.Range("F19").Value = -4
$ws.Range("F21").Value = -2
$ws.Range("F22").Value = -2

$wb.Save()
